$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J (shifts J:N -> K:O)
$ws.Columns("J").Insert()

# Insert a new row before row 11 (shifts old row 11 -> row 12)
$ws.Rows("11").Insert()

# Updated optimal solution value for E-n30-k3 (row 4)
$ws.Range("B4").Value = 534

# New row of results for X-n106-k14 (row 11)
$ws.Range("A11").Value = "X-n106-k14"
$ws.Range("B11").Value = 26362
$ws.Range("C11").Value = 14
$ws.Range("D11").Value = "25186,2…"
$ws.Range("E11").Value = "139s"

# Fix instance name in A2 (remove stray leading non-breaking space)
$ws.Range("A2").Value = "E-n22-k4"

# New ESPPRC result + runtime for E-n23-k3 (row 3)
$ws.Range("H3").Value = "558,9…"
$ws.Range("I3").Value = "102s"

# New header for the inserted column, and the comment for row 3
$ws.Range("J1").Value = "Kommentar"
$ws.Range("J3").Value = "time_limit=10, heuristic=10"

# Update selection to match the saved workbook state
$ws.Range("B5").Select()
